# Reposition the process-flow diagram picture on slide 3.
#
# Target OOXML (per the authoritative diff):
#   <a:off x="828675" y="2592019"/>  ->  <a:off x="1767500" y="2952369"/>
#
# PowerPoint's COM object model works in points (1 pt = 12700 EMU), so the
# new EMU offsets are converted to points below.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$pic = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Type -eq 13) {   # msoPicture
        $pic = $sh
        break
    }
}

# New offsets, in EMU, converted to points (EMU / 12700).
$pic.Left = 1767500 / 12700
$pic.Top  = 2952369 / 12700
